$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44223
$ws.Range("M2").Value = 200
$ws.Range("R2").Value = 'Provincia de Linares'

# Row 3
$ws.Range("D3").Value = 44196
$ws.Range("M3").Value = 150
$ws.Range("R3").Value = 'Provincia de Linares'

# Row 4
$ws.Range("D4").Value = 44193
$ws.Range("M4").Value = 200

# Row 5
$ws.Range("D5").Value = 44187
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 3400
$ws.Range("O5").Value = 3400
$ws.Range("P5").Value = 3400
$ws.Range("S5").Value = 1700

# Row 6
$ws.Range("D6").Value = 44187
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 4000
$ws.Range("P6").Value = 4000
$ws.Range("R6").Value = 'Provincia de Linares'
$ws.Range("S6").Value = 2000

# Row 7
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 3000
$ws.Range("O7").Value = 3000
$ws.Range("P7").Value = 3000
$ws.Range("R7").Value = 'Provincia de Linares'
$ws.Range("S7").Value = 1500

# Row 8
$ws.Range("D8").Value = 44253
$ws.Range("M8").Value = 25
$ws.Range("R8").Value = 'Provincia de Curicó'

# Row 9
$ws.Range("D9").Value = 44215
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 750
$ws.Range("N9").Value = 4000
$ws.Range("O9").Value = 4000
$ws.Range("P9").Value = 4000
$ws.Range("R9").Value = 'Provincia de Curicó'
$ws.Range("S9").Value = 2000

# Row 10
$ws.Range("D10").Value = 44188
$ws.Range("M10").Value = 300
$ws.Range("R10").Value = 'Provincia de Curicó'

# Row 11
$ws.Range("D11").Value = 44188
$ws.Range("M11").Value = 500
$ws.Range("N11").Value = 4000
$ws.Range("O11").Value = 4000
$ws.Range("P11").Value = 4000
$ws.Range("S11").Value = 2000

# Row 12
$ws.Range("D12").Value = 44221

# Row 13
$ws.Range("D13").Value = 44221

# Row 14
$ws.Range("D14").Value = 44224
$ws.Range("M14").Value = 250

# Row 15
$ws.Range("D15").Value = 44224
$ws.Range("M15").Value = 300
$ws.Range("R15").Value = 'Provincia de Linares'

# Row 16
$ws.Range("D16").Value = 44186

# Row 17
$ws.Range("D17").Value = 44250
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 4000
$ws.Range("P17").Value = 4000
$ws.Range("R17").Value = 'Provincia de Curicó'
$ws.Range("S17").Value = 2000

# Row 18
$ws.Range("D18").Value = 44175
$ws.Range("R18").Value = 'Provincia de Linares'

# Row 19
$ws.Range("D19").Value = 44203
$ws.Range("M19").Value = 350

# Row 20
$ws.Range("D20").Value = 44217
$ws.Range("M20").Value = 250

# Row 21
$ws.Range("D21").Value = 44217
$ws.Range("M21").Value = 300

# Row 22
$ws.Range("D22").Value = 44202
$ws.Range("M22").Value = 200
$ws.Range("N22").Value = 4000
$ws.Range("P22").Value = 4000
$ws.Range("S22").Value = 2000

# Row 23
$ws.Range("D23").Value = 44252
$ws.Range("M23").Value = 75
$ws.Range("R23").Value = 'Provincia de Curicó'

# Row 24
$ws.Range("D24").Value = 44211
$ws.Range("N24").Value = 3000
$ws.Range("O24").Value = 3500
$ws.Range("P24").Value = 3250
$ws.Range("R24").Value = 'Provincia de Curicó'
$ws.Range("S24").Value = 1625

# Row 25
$ws.Range("D25").Value = 44204
$ws.Range("M25").Value = 150
$ws.Range("N25").Value = 4000
$ws.Range("O25").Value = 4000
$ws.Range("P25").Value = 4000
$ws.Range("R25").Value = 'Provincia de Curicó'
$ws.Range("S25").Value = 2000

# Row 26
$ws.Range("D26").Value = 44204
$ws.Range("M26").Value = 250
$ws.Range("R26").Value = 'Provincia de Linares'

# Row 27
$ws.Range("D27").Value = 44260
$ws.Range("M27").Value = 75
$ws.Range("R27").Value = 'Provincia de Curicó'

# Row 28
$ws.Range("D28").Value = 44189
$ws.Range("M28").Value = 300
$ws.Range("N28").Value = 3000
$ws.Range("O28").Value = 3000
$ws.Range("P28").Value = 3000
$ws.Range("S28").Value = 1500

# Row 29
$ws.Range("D29").Value = 44189
$ws.Range("M29").Value = 250
$ws.Range("N29").Value = 3000
$ws.Range("O29").Value = 3000
$ws.Range("P29").Value = 3000
$ws.Range("S29").Value = 1500

# Row 30
$ws.Range("D30").Value = 44209
$ws.Range("M30").Value = 170
$ws.Range("N30").Value = 3000
$ws.Range("P30").Value = 3500
$ws.Range("R30").Value = 'Provincia de Linares'
$ws.Range("S30").Value = 1750

# Row 31
$ws.Range("D31").Value = 44216
$ws.Range("M31").Value = 200
$ws.Range("R31").Value = 'Provincia de Curicó'

# Row 32
$ws.Range("D32").Value = 44216
$ws.Range("M32").Value = 400
$ws.Range("N32").Value = 4000
$ws.Range("O32").Value = 4000
$ws.Range("P32").Value = 4000
$ws.Range("R32").Value = 'Provincia de Linares'
$ws.Range("S32").Value = 2000

# Row 33
$ws.Range("D33").Value = 44257
$ws.Range("M33").Value = 100
$ws.Range("N33").Value = 4000
$ws.Range("O33").Value = 4000
$ws.Range("P33").Value = 4000
$ws.Range("R33").Value = 'Provincia de Curicó'
$ws.Range("S33").Value = 2000

# Row 34
$ws.Range("D34").Value = 44169
$ws.Range("M34").Value = 200
$ws.Range("N34").Value = 5000
$ws.Range("O34").Value = 5000
$ws.Range("P34").Value = 5000
$ws.Range("S34").Value = 2500

# Row 35
$ws.Range("D35").Value = 44195
$ws.Range("M35").Value = 300
$ws.Range("N35").Value = 3000
$ws.Range("O35").Value = 3000
$ws.Range("P35").Value = 3000
$ws.Range("S35").Value = 1500

# Row 36
$ws.Range("D36").Value = 44239
$ws.Range("M36").Value = 350
$ws.Range("N36").Value = 3500
$ws.Range("P36").Value = 3750
$ws.Range("S36").Value = 1875

# Row 37
$ws.Range("D37").Value = 44222
$ws.Range("R37").Value = 'Provincia de Curicó'

# Row 38
$ws.Range("D38").Value = 44222
$ws.Range("M38").Value = 300
$ws.Range("R38").Value = 'Provincia de Linares'

# Row 39
$ws.Range("D39").Value = 44210
$ws.Range("M39").Value = 400
$ws.Range("N39").Value = 3000
$ws.Range("P39").Value = 3500
$ws.Range("S39").Value = 1750

# Row 40
$ws.Range("D40").Value = 44176
$ws.Range("M40").Value = 100
$ws.Range("N40").Value = 4000
$ws.Range("O40").Value = 4000
$ws.Range("P40").Value = 4000
$ws.Range("R40").Value = 'Provincia de Linares'
$ws.Range("S40").Value = 2000

# Row 41
$ws.Range("D41").Value = 44225
$ws.Range("M41").Value = 150

# Row 42
$ws.Range("D42").Value = 44225

# Row 43
$ws.Range("D43").Value = 44251
$ws.Range("M43").Value = 125
$ws.Range("R43").Value = 'Provincia de Curicó'

# Row 45
$ws.Range("D45").Value = 44194
$ws.Range("R45").Value = 'Provincia de Linares'

# Row 46
$ws.Range("D46").Value = 44201
$ws.Range("M46").Value = 200
